$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: add a new paragraph after "Unlike real life Floor Control..."
#   "System needs updated concurrency changes to handle multiple elevators" +
#   bookmark _GoBack + "."
# ---------------------------------------------------------------------------
$oldPara1 = "Unlike real life Floor Control doesn’t limit UP and DOWN button on top and ground floor respectively"
$newText1 = "System needs updated concurrency changes to handle multiple elevators."
$rng1 = $d.Content
$ok1 = $rng1.Find.Execute($oldPara1, $true, $false, $false, $false, $false, $true, 1, $false, ($oldPara1 + "^p" + $newText1), 2)

# Locate the freshly inserted paragraph and place the _GoBack bookmark right
# before the trailing period (between "elevators" and ".").
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -eq ($newText1 + "`r")) {
        $newPara1 = $d.Paragraphs($i)
        break
    }
}
$bmPos = $newPara1.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# Change 2: "Possible improvements in design and implementation" ->
#   "Possible improvements" + " in " + "design, implementation" + " and testing"
#   (each new run keeps the existing bold formatting, split into 3 runs)
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -eq "Possible improvements in design and implementation`r") {
        $para2 = $d.Paragraphs($i)
        break
    }
}
$marker = "Possible improvements"
$pStart = $para2.Range.Start
$runStart = $pStart + $marker.Length
$runEnd = $para2.Range.End - 1

# Turn bold off before replacing the text so the new text does not get
# coalesced into the preceding "Possible improvements" run.
$d.Range($runStart, $runEnd).Font.Bold = 0
$d.Range($runStart, $runEnd).Text = " in design, implementation and testing"

$para2 = $d.Paragraphs($i)
$pStart = $para2.Range.Start
$s1 = $pStart + $marker.Length
$e1 = $s1 + " in ".Length
$s2 = $e1
$e2 = $s2 + "design, implementation".Length
$s3 = $e2
$e3 = $para2.Range.End - 1

# Re-apply bold to the whole replaced span ...
$d.Range($s1, $e3).Font.Bold = 1
# ... then toggle the middle segment off/on so it stays a distinct run
# instead of being re-merged with its (identically formatted) neighbours.
$d.Range($s2, $e2).Font.Bold = 0
$d.Range($s2, $e2).Font.Bold = 1

# ---------------------------------------------------------------------------
# Change 3: add a new paragraph after "...and others that supports observer
#   testing" -> "Test cases for concurrency management can be added."
# ---------------------------------------------------------------------------
$oldPara3 = "and others that supports observer testing"
$newText3 = "Test cases for concurrency management can be added."
$rng3 = $d.Content
$ok3 = $rng3.Find.Execute($oldPara3, $true, $false, $false, $false, $false, $true, 1, $false, ($oldPara3 + "^p" + $newText3), 2)

# ---------------------------------------------------------------------------
# Change 4: normalise the "currentFloor=3, upQueue=..." run split so the
#   "=3" and ", " runs are merged into a single "=3, " run (matching the
#   canonical OOXML produced once the trailing _GoBack bookmark moved away).
# ---------------------------------------------------------------------------
$rng4 = $d.Content
$ok4 = $rng4.Find.Execute("=3", $true, $false, $false, $false, $false, $true, 1, $false, "=3", 2)

Write-Host "Change1: $ok1  Change3: $ok3  Change4: $ok4"
